$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 563.2727
$ws.Range("I2").Value = 376.7143
$ws.Range("J2").Value = 889.75
$ws.Range("K2").Value = 376.7143
$ws.Range("L2").Value = 889.75
$ws.Range("M2").Value = -263.7143
$ws.Range("N2").Value = -1115.75
$ws.Range("H19").Value = 1849.5714
$ws.Range("I19").Value = 349.83334
$ws.Range("K19").Value = 349.83334
$ws.Range("M19").Value = -174.83334
$ws.Range("H41").Value = 350.125
$ws.Range("I41").Value = 371.57144
$ws.Range("K41").Value = 371.57144
$ws.Range("M41").Value = 68.42856
$ws.Range("H51").Value = 53984.477
$ws.Range("I51").Value = 10199
$ws.Range("J51").Value = 71498.664
$ws.Range("K51").Value = 10199
$ws.Range("L51").Value = 71498.664
$ws.Range("M51").Value = -9715
$ws.Range("N51").Value = -72466.664
$ws.Range("H98").Value = 2443.4285
$ws.Range("I98").Value = 2443.4285
$ws.Range("K98").Value = 2443.4285
$ws.Range("M98").Value = -945.4285
$ws.Range("J111").Value = 2600
$ws.Range("L111").Value = 7800
$ws.Range("N111").Value = -13934
$ws.Range("H122").Value = 2443.4285
$ws.Range("I122").Value = 2443.4285
$ws.Range("K122").Value = 7330.2855
$ws.Range("M122").Value = -4880.2855
$ws.Range("H131").Value = 624.75
$ws.Range("I131").Value = 624.75
$ws.Range("K131").Value = 1874.25
$ws.Range("M131").Value = 3165.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17545110
$ws.Range("I2").Value = 20834076
$ws.Range("K2").Value = 20834076
$ws.Range("M2").Value = -20833963
$ws.Range("H5").Value = 1931.6666
$ws.Range("I5").Value = 1897.5
$ws.Range("K5").Value = 1897.5
$ws.Range("M5").Value = -1785.5
$ws.Range("H10").Value = 5264
$ws.Range("I10").Value = 3962
$ws.Range("K10").Value = 3962
$ws.Range("M10").Value = -3792
$ws.Range("H12").Value = 1625000
$ws.Range("I12").Value = 1625000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1625000
$ws.Range("L12").ClearContents()
$ws.Range("N12").Value = 0
$ws.Range("M12").Value = -1624827
$ws.Range("H110").Value = 3970580.5
$ws.Range("I110").Value = 6537462
$ws.Range("K110").Value = 6537462
$ws.Range("M110").Value = -6535417
$ws.Range("H116").Value = 17545110
$ws.Range("I116").Value = 20834076
$ws.Range("K116").Value = 20834076
$ws.Range("M116").Value = -20831782
$ws.Range("H132").Value = 3634.1765
$ws.Range("I132").Value = 4070.7856
$ws.Range("K132").Value = 12212.3568
$ws.Range("M132").Value = -9682.356800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17545110
$ws.Range("I3").Value = 20834076
$ws.Range("K3").Value = 20834076
$ws.Range("M3").Value = -20833962
$ws.Range("H4").Value = 1931.6666
$ws.Range("I4").Value = 1897.5
$ws.Range("K4").Value = 1897.5
$ws.Range("M4").Value = -1782.5
$ws.Range("H5").Value = 1248.5
$ws.Range("I5").Value = 165
$ws.Range("J5").Value = 4499
$ws.Range("K5").Value = 165
$ws.Range("L5").Value = 4499
$ws.Range("M5").Value = -52
$ws.Range("N5").Value = -4725
$ws.Range("H86").Value = 1389.8
$ws.Range("I86").Value = 1384.4286
$ws.Range("K86").Value = 1384.4286
$ws.Range("M86").Value = -261.4286
$ws.Range("H89").Value = 1389.8
$ws.Range("I89").Value = 1384.4286
$ws.Range("K89").Value = 6922.143
$ws.Range("M89").Value = -1306.143
$ws.Range("H94").Value = 2521.7693
$ws.Range("I94").Value = 1473.75
$ws.Range("K94").Value = 1473.75
$ws.Range("M94").Value = -1022.75
$ws.Range("H105").Value = 3791826
$ws.Range("I105").Value = 5558864.5
$ws.Range("K105").Value = 5558864.5
$ws.Range("M105").Value = -5557117.5
$ws.Range("H134").Value = 2017.6428
$ws.Range("I134").Value = 2186.2727
$ws.Range("K134").Value = 6558.8181
$ws.Range("M134").Value = -4023.8181

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 550000
$ws.Range("I22").Value = 999999
$ws.Range("J22").Value = 100001
$ws.Range("K22").Value = 999999
$ws.Range("L22").Value = 100001
$ws.Range("M22").Value = -999649
$ws.Range("N22").Value = -100701
$ws.Range("H31").Value = 3511.8333
$ws.Range("I31").Value = 2692.1738
$ws.Range("K31").Value = 2692.1738
$ws.Range("M31").Value = -2397.1738
$ws.Range("H34").Value = 3511.8333
$ws.Range("I34").Value = 2692.1738
$ws.Range("K34").Value = 2692.1738
$ws.Range("M34").Value = -2490.1738
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H62").Value = 103749.75
$ws.Range("H65").Value = 103749.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 16666786
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H39").Value = 9632.736999999999
$ws.Range("J39").Value = 12270.923
$ws.Range("L39").Value = 36812.769
$ws.Range("N39").Value = -37400.769
$ws.Range("H55").Value = 9583.333000000001
$ws.Range("J55").Value = 10187.5
$ws.Range("L55").Value = 30562.5
$ws.Range("N55").Value = -30916.5
$ws.Range("H88").Value = 5069.5
$ws.Range("I88").Value = 5069.5
$ws.Range("K88").Value = 15208.5
$ws.Range("M88").Value = -14780.5
$ws.Range("H91").Value = 5069.5
$ws.Range("I91").Value = 5069.5
$ws.Range("K91").Value = 15208.5
$ws.Range("M91").Value = -13726.5
$ws.Range("H92").Value = 533.3333
$ws.Range("J92").Value = 500
$ws.Range("L92").Value = 1500
$ws.Range("N92").Value = -3996
$ws.Range("H109").Value = 1294.5
$ws.Range("I109").Value = 953.6
$ws.Range("J109").Value = 2999
$ws.Range("K109").Value = 2860.8
$ws.Range("L109").Value = 8997
$ws.Range("M109").Value = -1820.8
$ws.Range("N109").Value = -11077

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 7579.75
$ws.Range("J7").Value = 7579.75
$ws.Range("L7").Value = 7579.75
$ws.Range("N7").Value = -7803.75
$ws.Range("H8").Value = 7579.75
$ws.Range("J8").Value = 7579.75
$ws.Range("L8").Value = 7579.75
$ws.Range("N8").Value = -7857.75
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -860
$ws.Range("H13").Value = 150
$ws.Range("I13").Value = 150
$ws.Range("K13").Value = 150
$ws.Range("M13").Value = -11
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 10000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9832
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H41").Value = 4000
$ws.Range("I41").Value = 4000
$ws.Range("K41").Value = 4000
$ws.Range("M41").Value = -3645
$ws.Range("H126").Value = 8099
$ws.Range("I126").Value = 7748.75
$ws.Range("J126").Value = 9500
$ws.Range("K126").Value = 23246.25
$ws.Range("L126").Value = 28500
$ws.Range("M126").Value = -20776.25
$ws.Range("N126").Value = -33440
$ws.Range("H132").Value = 4997.8
$ws.Range("I132").Value = 4997.8
$ws.Range("K132").Value = 14993.4
$ws.Range("M132").Value = -12463.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 168888.83
$ws.Range("I132").Value = 202222.2
$ws.Range("K132").Value = 606666.6000000001
$ws.Range("M132").Value = -604136.6000000001
$ws.Range("H136").Value = 5821.077
$ws.Range("I136").Value = 4834.1113
$ws.Range("K136").Value = 14502.3339
$ws.Range("M136").Value = -11952.3339

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 667033.3
$ws.Range("H4").Value = 4200300
$ws.Range("H113").Value = 1237.3572
$ws.Range("I113").Value = 391.14285
$ws.Range("J113").Value = 2083.5715
$ws.Range("K113").Value = 1173.42855
$ws.Range("L113").Value = 6250.7145
$ws.Range("M113").Value = 996.5714499999999
$ws.Range("N113").Value = -10590.7145
$ws.Range("H126").Value = 834.3333
$ws.Range("I126").Value = 834.3333
$ws.Range("K126").Value = 2502.9999
$ws.Range("M126").Value = -32.9998999999998
$ws.Range("H132").Value = 1964.2
$ws.Range("I132").Value = 1888.625
$ws.Range("J132").Value = 2266.5
$ws.Range("K132").Value = 5665.875
$ws.Range("L132").Value = 6799.5
$ws.Range("M132").Value = -3135.875
$ws.Range("N132").Value = -11859.5
$ws.Range("H136").Value = 2027.5555
$ws.Range("I136").Value = 1026
$ws.Range("J136").Value = 4030.6667
$ws.Range("K136").Value = 3078
$ws.Range("L136").Value = 12092.0001
$ws.Range("M136").Value = -528
$ws.Range("N136").Value = -17192.0001
